# Edit: update "as of" date in confidential disclosure note, and refresh the
# Weight / Percent Change model values for each holding (rows 2-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (with a password) - unprotect for the duration of the
# edit so the cell values can be written, then restore protection afterwards.
$ws.Unprotect("D382")

# Update the "Model holdings provided as of ..." confidential disclosure note.
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# New Weight (column D) / Percent Change (column E) values for rows 2-41.
$holdingUpdates = @(
    @{r=2; d=0.07897453946587876; e=0.01923377638780299},
    @{r=3; d=0.06980509528238932; e=0.01340536214485777},
    @{r=4; d=0.05425923402135772; e=-0.0004172672786587617},
    @{r=5; d=0.04821257585520834; e=0.006071251055836724},
    @{r=6; d=0.04327692441816958; e=0.00122636029174461},
    @{r=7; d=0.03868367534610351; e=0.006521639987230099},
    @{r=8; d=0.03848549656208325; e=-0.00391174133610428},
    @{r=9; d=0.03494003695122817; e=-0.0006437768240343811},
    @{r=10; d=0.03291752542914601; e=0.005091490511516028},
    @{r=11; d=0.02757449015640981; e=-0.001279590531030106},
    @{r=12; d=0.03179534541547637; e=-0.008739076154806291},
    @{r=13; d=0.03222066531322906; e=-0.0002235778849530856},
    @{r=14; d=0.02727751601418359; e=0.01794761237469022},
    @{r=15; d=0.03015463678812762; e=0.003802837501827971},
    @{r=16; d=0.02687763004049286; e=-0.002051197899573309},
    @{r=17; d=0.02756302283804662; e=0.006272602169808073},
    @{r=18; d=0.02337568743265707; e=0.01207547169811307},
    @{r=19; d=0.01972878615935838; e=0.01912157026394601},
    @{r=20; d=0.02167323170642959; e=-0.004945054945055039},
    @{r=21; d=0.02079921622349156; e=-0.01042587029510511},
    @{r=22; d=0.02168499305859696; e=-0.02372881355932199},
    @{r=23; d=0.0203659574130259; e=-0.003003003003003046},
    @{r=24; d=0.02000782423952934; e=-0.03006789524733278},
    @{r=25; d=0.01765820011029208; e=-0.002747481475314228},
    @{r=26; d=0.01806631903049997; e=-0.02094624285923519},
    @{r=27; d=0.01902766255328076; e=0.002704268881591698},
    @{r=28; d=0.01684946013188298; e=-0.01122948459545048},
    @{r=29; d=0.01790548253961113; e=-0.008276405675249787},
    @{r=30; d=0.01742194394862994; e=0.001392369813422611},
    @{r=31; d=0.01851222129454556; e=0.01633589847441619},
    @{r=32; d=0.01546764826911855; e=0.008459271932325985},
    @{r=33; d=0.01686695514323195; e=-0.006066522557701681},
    @{r=34; d=0.008317334218963151; e=0.01226711917135059},
    @{r=35; d=0.008041677527540307; e=0.01387593923106456},
    @{r=36; d=0.007515504034952385; e=0.03482003129890465},
    @{r=37; d=0.006491972362586608; e=0.01657683771909957},
    @{r=38; d=0.007092977458339451; e=0.006259586286946117},
    @{r=39; d=0.007253960966130392; e=0.01303175857805861},
    @{r=40; d=0.006856574279775221; e=0.01099961404862992},
    @{r=41; d=0.9999999999999998; e=0.003024137676035599}
)

foreach ($item in $holdingUpdates) {
    $ws.Cells.Item($item.r, 4).Value = $item.d
    $ws.Cells.Item($item.r, 5).Value = $item.e
}

# Restore sheet protection with the original password.
$ws.Protect("D382")
